# Normalize the "Recorded By" (column G) entries on the
# "Session Analysis Results" sheet: for every multi-author cell whose
# last listed author is not already "System", swap the order of the
# final two comma-separated names/emails (effectively moving the
# recorder that used to be listed last, in front of "System").
#
# Example:
#   "System, dnasr281@gmail.com"              -> "dnasr281@gmail.com, System"
#   "system, System, backup@backdoor.com"      -> "system, backup@backdoor.com, System"
#   "admin@admin.com, dnasr281@gmail.com"      -> "dnasr281@gmail.com, admin@admin.com"
#   "backup@backdoor.com, System"              -> unchanged (already ends with "System")

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)   # Column G = "Recorded By"
    $value = $cell.Value2

    if ([string]::IsNullOrEmpty($value)) {
        continue
    }

    $parts = $value -split ",\s*"

    if ($parts.Count -ge 2 -and $parts[$parts.Count - 1] -ne "System") {
        $last = $parts[$parts.Count - 1]
        $secondLast = $parts[$parts.Count - 2]
        $parts[$parts.Count - 1] = $secondLast
        $parts[$parts.Count - 2] = $last

        $newValue = [string]::Join(", ", $parts)
        $cell.Value2 = $newValue
    }
}
